# Apply updated simulation result values to the "variants" and "results" sheets.
$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    # Force the cell to hold a text value (matching original inlineStr / shared-string
    # cells) even when the text looks numeric, then strip the formatting change that
    # switching the number format introduced so the cell style is left untouched.
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.ClearFormats()
}

# --- "variants" sheet (numeric cells, columns A:C, style index 3 preserved) ---
$wsVariants = $wb.Worksheets.Item("variants")
$wsVariants.Cells.Item(2, 2).Value2 = 5.165942130649757
$wsVariants.Cells.Item(2, 3).Value2 = 11154.18953864742
$wsVariants.Cells.Item(3, 2).Value2 = 4.798849819416276
$wsVariants.Cells.Item(3, 3).Value2 = 13672.96917189622

# --- "results" sheet (columns A:I; B/C numeric, D:I text) ---
$wsResults = $wb.Worksheets.Item("results")

# Row 2
$wsResults.Cells.Item(2, 2).Value2 = 5.165942130649757
$wsResults.Cells.Item(2, 3).Value2 = 11154.18953864742
Set-TextValue $wsResults.Cells.Item(2, 4) "276"
Set-TextValue $wsResults.Cells.Item(2, 5) "11153.7"
Set-TextValue $wsResults.Cells.Item(2, 6) "11157.6"
Set-TextValue $wsResults.Cells.Item(2, 7) "276"
Set-TextValue $wsResults.Cells.Item(2, 8) "3.45664"
Set-TextValue $wsResults.Cells.Item(2, 9) "3.12996"

# Row 3
$wsResults.Cells.Item(3, 2).Value2 = 4.798849819416276
$wsResults.Cells.Item(3, 3).Value2 = 13672.96917189622
Set-TextValue $wsResults.Cells.Item(3, 4) "275"
Set-TextValue $wsResults.Cells.Item(3, 5) "13672.6"
Set-TextValue $wsResults.Cells.Item(3, 6) "13676"
Set-TextValue $wsResults.Cells.Item(3, 7) "275"
Set-TextValue $wsResults.Cells.Item(3, 8) "3.22726"
Set-TextValue $wsResults.Cells.Item(3, 9) "2.91238"
